# The deck's only table lives on slide 5 (the B1 "Types of financial
# documents" slide) as the 2nd shape on that slide (Shape 1 is the title
# textbox, Shape 2 is the graphicFrame/table, Shape 3 is another textbox).
#
# The author picked a different built-in table style for it via the Table
# Design ribbon, which changes <a:tblPr>/<a:tableStyleId> from the old
# GUID to the new one. In the PowerPoint object model that is expressed
# through Table.ApplyStyle(StyleId, bandedRows).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(5)
$shape = $s.Shapes.Item(2)
$tbl = $shape.Table

$tbl.ApplyStyle("{6F28097B-6518-4EE5-AC82-A043CD66E1EB}", $true)
